# Apply updated "想去人数" (F column) counts to the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Row => new F-column value for worksheet "展览" (sheet1)
$changesExhibition = @{
    3  = 556
    5  = 112
    8  = 50
    10 = 16148
    12 = 193
    14 = 6307
    29 = 42
    30 = 5036
    32 = 11260
    33 = 1243
    35 = 141
}

# Row => new F-column value for worksheet "全部类型" (sheet4)
$changesAllTypes = @{
    3  = 556
    5  = 112
    8  = 50
    10 = 16148
    12 = 193
    14 = 6307
    29 = 42
    30 = 5036
    33 = 11260
    34 = 1243
    36 = 141
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $changesExhibition.Keys) {
    $wsExhibition.Range("F$row").Value = $changesExhibition[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $changesAllTypes.Keys) {
    $wsAllTypes.Range("F$row").Value = $changesAllTypes[$row]
}
